{"js": "// Change 1 (4 occurrences): drop the redundant \"constel\u00b7laci\u00f3, \" lead-in\n// and lower-case \"Bessons\" -> \"bessons\" in the \"Dates de la campanya...\"\n// sentence.\nconst oldA = \"la constel\u00b7laci\u00f3, Constel\u00b7laci\u00f3 de Bessons\";\nconst newA = \"la  Constel\u00b7laci\u00f3 de bessons\";\n\nconst hitsA = context.document.body.search(oldA, { matchCase: true });\nhitsA.load(\"items\");\nawait context.sync();\n\nfor (const r of hitsA.items) {\n  r.insertText(newA, \"Replace\");\n}\nawait context.sync();\n\n// Change 2 (1 occurrence): lower-case \"Bessons\" -> \"bessons\" inside the\n// longer \"Esteu participant...\" paragraph.\nconst oldB = \"Constel\u00b7laci\u00f3 de Bessons a la nit\";\nconst newB = \"Constel\u00b7laci\u00f3 de bessons a la nit\";\n\nconst hitsB = context.document.body.search(oldB, { matchCase: true });\nhitsB.load(\"items\");\nawait context.sync();\n\nfor (const r of hitsB.items) {\n  r.insertText(newB, \"Replace\");\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Change 1 (4 occurrences): drop the redundant \"constel\u00b7laci\u00f3, \" lead-in\n# and lower-case \"Bessons\" -> \"bessons\" in the \"Dates de la campanya...\"\n# sentence.\n$rng1 = $d.Content\n$find1 = $rng1.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"la constel\u00b7laci\u00f3, Constel\u00b7laci\u00f3 de Bessons\"\n$find1.Replacement.Text = \"la  Constel\u00b7laci\u00f3 de bessons\"\n$find1.Forward = $true\n$find1.Wrap = 1\n$find1.Format = $false\n$find1.MatchCase = $true\n$find1.MatchWholeWord = $false\n$find1.MatchWildcards = $false\n[void]$find1.Execute($find1.Text, $find1.MatchCase, $find1.MatchWholeWord, $find1.MatchWildcards, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n# Change 2 (1 occurrence): lower-case \"Bessons\" -> \"bessons\" inside the\n# longer \"Esteu participant...\" paragraph.\n$rng2 = $d.Content\n$find2 = $rng2.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"Constel\u00b7laci\u00f3 de Bessons a la nit\"\n$find2.Replacement.Text = \"Constel\u00b7laci\u00f3 de bessons a la nit\"\n$find2.Forward = $true\n$find2.Wrap = 1\n$find2.Format = $false\n$find2.MatchCase = $true\n$find2.MatchWholeWord = $false\n$find2.MatchWildcards = $false\n[void]$find2.Execute($find2.Text, $find2.MatchCase, $find2.MatchWholeWord, $find2.MatchWildcards, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
